$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended below the existing table (row 25).
# Column A holds the date as plain text (matching the format used by the
# other rows in this column), so force a Text number format before
# assigning the value to stop Excel auto-converting it to a date serial,
# then restore the cell's style back to Normal/default.
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "01/30/2026"
$ws.Range("A25").Style = "Normal"

$ws.Range("B25").Value = 1275.224000000002
$ws.Range("C25").Value = 0.03881670984862261
$ws.Range("D25").Value = 50
